# chore: update Sheets via scheduled runner
#
# Refreshes cached Universalis market-price figures (currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ) and the resulting leve-profit
# figures (LeveProfitNQ/HQ) for a batch of leves across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets of Masamune_Profits. These are plain
# cached values (no formulas in this workbook), so each refreshed leve is
# applied as a direct cell write per sheet/row.

# Auto-generated Excel COM-interop script to apply the Masamune_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 95: Official Strategy Guide
$ws.Range("H95").Value = 34657.332
$ws.Range("J95").Value = 34657.332
$ws.Range("L95").Value = 34657.332
$ws.Range("N95").Value = -40149.332

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 19105.283
$ws.Range("I132").Value = 2833.1226
$ws.Range("J132").Value = 218439.25
$ws.Range("K132").Value = 8499.3678
$ws.Range("L132").Value = 655317.75
$ws.Range("M132").Value = -5969.3678
$ws.Range("N132").Value = -660377.75

# Row 134: Binding Spells
$ws.Range("H134").Value = 81100
$ws.Range("J134").Value = 81100
$ws.Range("L134").Value = 81100
$ws.Range("N134").Value = -91240

# Row 135: For Tired Minds
$ws.Range("H135").Value = 11628797
$ws.Range("I135").Value = 904.44446
$ws.Range("K135").Value = 8140.00014
$ws.Range("M135").Value = -5605.00014

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 3269.32
$ws.Range("I137").Value = 1340.75
$ws.Range("J137").Value = 3499.597
$ws.Range("K137").Value = 4022.25
$ws.Range("L137").Value = 10498.791
$ws.Range("M137").Value = -1472.25
$ws.Range("N137").Value = -15598.791

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3049.082
$ws.Range("I138").Value = 2236.6538
$ws.Range("J138").Value = 3652.6
$ws.Range("K138").Value = 6709.9614
$ws.Range("L138").Value = 10957.8
$ws.Range("M138").Value = -1569.9614
$ws.Range("N138").Value = -21237.8

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 18034.13
$ws.Range("I32").Value = 16471.725
$ws.Range("J32").Value = 26272.273
$ws.Range("K32").Value = 16471.725
$ws.Range("L32").Value = 26272.273
$ws.Range("M32").Value = -16184.725
$ws.Range("N32").Value = -26846.273

# Row 95: Shielded Life
$ws.Range("H95").Value = 37195
$ws.Range("J95").Value = 37195
$ws.Range("L95").Value = 37195
$ws.Range("N95").Value = -42687

# Row 96: The Gauntlet Is Cast
$ws.Range("H96").Value = 29125.8
$ws.Range("J96").Value = 29125.8
$ws.Range("L96").Value = 29125.8
$ws.Range("N96").Value = -34617.8

# Row 103: Sweeping the Legs
$ws.Range("H103").Value = 38300
$ws.Range("J103").Value = 38300
$ws.Range("L103").Value = 38300
$ws.Range("N103").Value = -40644

# Row 104: See Shields by the Sea Shore
$ws.Range("H104").Value = 29369.334
$ws.Range("J104").Value = 29369.334
$ws.Range("L104").Value = 29369.334
$ws.Range("N104").Value = -36357.334

# Row 106: Heads Will Roll
$ws.Range("H106").Value = 39996
$ws.Range("J106").Value = 39996
$ws.Range("L106").Value = 39996
$ws.Range("N106").Value = -42520

# Row 135: Forgiveness for My Shins
$ws.Range("H135").Value = 49189.57
$ws.Range("J135").Value = 49189.57
$ws.Range("L135").Value = 49189.57
$ws.Range("N135").Value = -59329.57

$ws = $wb.Worksheets.Item("BSM")
# Row 92: Have Blade, Will Travel
$ws.Range("H92").Value = 38650
$ws.Range("J92").Value = 38650
$ws.Range("L92").Value = 38650
$ws.Range("N92").Value = -43642

# Row 100: And My Axe
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2273.6562
$ws.Range("I105").Value = 2269
$ws.Range("J105").Value = 2277.7646
$ws.Range("K105").Value = 2269
$ws.Range("L105").Value = 2277.7646
$ws.Range("M105").Value = -522
$ws.Range("N105").Value = -5771.7646

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 15389.667
$ws.Range("I31").Value = 4749.1665
$ws.Range("J31").Value = 36670.668
$ws.Range("K31").Value = 4749.1665
$ws.Range("L31").Value = 36670.668
$ws.Range("M31").Value = -4454.1665
$ws.Range("N31").Value = -37260.668

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 15389.667
$ws.Range("I34").Value = 4749.1665
$ws.Range("J34").Value = 36670.668
$ws.Range("K34").Value = 4749.1665
$ws.Range("L34").Value = 36670.668
$ws.Range("M34").Value = -4547.1665
$ws.Range("N34").Value = -37074.668

# Row 43: The Long Lance of the Law
$ws.Range("H43").Value = 49657
$ws.Range("J43").Value = 49657
$ws.Range("L43").Value = 49657
$ws.Range("N43").Value = -50025

# Row 92: Walk the Walk
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 96: Composition
$ws.Range("H96").Value = 93538.664
$ws.Range("J96").Value = 93538.664
$ws.Range("L96").Value = 93538.664
$ws.Range("N96").Value = -99030.664

# Row 101: Everybody's Heard about the 'Berd
$ws.Range("H101").Value = 49657
$ws.Range("J101").Value = 49657
$ws.Range("L101").Value = 49657
$ws.Range("N101").Value = -56147

$ws = $wb.Worksheets.Item("CUL")
# Row 70: Persona non Gratin
$ws.Range("H70").Value = 5666.6665
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5666.6665
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 16999.9995
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -17629.9995

# Row 73: Recipe for Disaster (L)
$ws.Range("H73").Value = 5666.6665
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5666.6665
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 16999.9995
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -19183.9995

# Row 88: Don't Let It Fall Apart
$ws.Range("H88").Value = 7014.933
$ws.Range("J88").Value = 7014.933
$ws.Range("L88").Value = 21044.799
$ws.Range("N88").Value = -21900.799

# Row 91: Better Come Back with a Sandwich (L)
$ws.Range("H91").Value = 7014.933
$ws.Range("J91").Value = 7014.933
$ws.Range("L91").Value = 21044.799
$ws.Range("N91").Value = -24008.799

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 51376.69
$ws.Range("I131").Value = 20585.8
$ws.Range("J131").Value = 55537.62
$ws.Range("K131").Value = 61757.39999999999
$ws.Range("L131").Value = 166612.86
$ws.Range("M131").Value = -56717.39999999999
$ws.Range("N131").Value = -176692.86

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 38468696
$ws.Range("I137").Value = 2566.75
$ws.Range("J137").Value = 100014504
$ws.Range("K137").Value = 7700.25
$ws.Range("L137").Value = 300043512
$ws.Range("M137").Value = -2600.25
$ws.Range("N137").Value = -300053712

$ws = $wb.Worksheets.Item("GSM")
# Row 96: Bracelet for Impact
$ws.Range("H96").Value = 37261
$ws.Range("J96").Value = 37261
$ws.Range("L96").Value = 37261
$ws.Range("N96").Value = -42753

# Row 105: Untucked
$ws.Range("H105").Value = 40131
$ws.Range("J105").Value = 40131
$ws.Range("L105").Value = 40131
$ws.Range("N105").Value = -47119

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 2172.65
$ws.Range("I7").Value = 1849.8667
$ws.Range("J7").Value = 3141
$ws.Range("K7").Value = 1849.8667
$ws.Range("L7").Value = 3141
$ws.Range("M7").Value = -1737.8667
$ws.Range("N7").Value = -3365

# Row 111: Glove Me Tender
$ws.Range("H111").Value = 38333.25
$ws.Range("J111").Value = 38333.25
$ws.Range("L111").Value = 38333.25
$ws.Range("N111").Value = -46513.25

# Row 126: Battered Books
$ws.Range("H126").Value = 2172.65
$ws.Range("I126").Value = 1849.8667
$ws.Range("J126").Value = 3141
$ws.Range("K126").Value = 5549.6001
$ws.Range("L126").Value = 9423
$ws.Range("M126").Value = -3079.6001
$ws.Range("N126").Value = -14363

$ws = $wb.Worksheets.Item("WVR")
# Row 94: Proper Props
$ws.Range("H94").Value = 23745
$ws.Range("J94").Value = 23745
$ws.Range("L94").Value = 23745
$ws.Range("N94").Value = -25547

